$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.779.15"
$ws.Range("E2").Value = "'  -0.95%  "
$ws.Range("D3").Value = "'1.626.63"
$ws.Range("E3").Value = "'  -0.91%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "'  +0.12%  "
$ws.Range("D5").Value = "'215.00"
$ws.Range("E5").Value = "'  -0.44%  "
$ws.Range("D6").Value = "'0.5111"
$ws.Range("E6").Value = "'  -0.92%  "
$ws.Range("E7").Value = "'  +0.09%  "
$ws.Range("D8").Value = "'0.2581"
$ws.Range("E8").Value = "'  -0.23%  "
$ws.Range("D9").Value = "'0.06376"
$ws.Range("D10").Value = "'19.34"
$ws.Range("E10").Value = "'  -2.75%  "
$ws.Range("D11").Value = "'0.07780"
$ws.Range("E11").Value = "'  +0.09%  "
$ws.Range("D12").Value = "'4.252"
$ws.Range("E12").Value = "'  -1.03%  "
$ws.Range("D13").Value = "'1.628.10"
$ws.Range("E13").Value = "'  -0.90%  "
$ws.Range("D14").Value = "'1.849.22"
$ws.Range("E14").Value = "'  -1.21%  "
$ws.Range("D15").Value = "'0.5553"
$ws.Range("E15").Value = "'  +1.30%  "
$ws.Range("D16").Value = "'63.50"
$ws.Range("E16").Value = "'  -1.74%  "
$ws.Range("D17").Value = "'0.0₅7532"
$ws.Range("E17").Value = "'  -3.41%  "
$ws.Range("D18").Value = "'25.796.26"
$ws.Range("E18").Value = "'  -0.91%  "
$ws.Range("D19").Value = "'1.004"
$ws.Range("E19").Value = "'  +0.12%  "
$ws.Range("D20").Value = "'193.79"
$ws.Range("E20").Value = "'  -2.94%  "
$ws.Range("D21").Value = "'4.333"
$ws.Range("E21").Value = "'  -3.18%  "
$ws.Range("D22").Value = "'9.786"
$ws.Range("E22").Value = "'  -2.13%  "
$ws.Range("D23").Value = "'5.993"
$ws.Range("E23").Value = "'  -2.00%  "
$ws.Range("E24").Value = "'  -0.04%  "
$ws.Range("D25").Value = "'1.820"
$ws.Range("E25").Value = "'  -4.38%  "
$ws.Range("D26").Value = "'0.1294"
$ws.Range("E26").Value = "'  +5.00%  "
$ws.Range("D27").Value = "'141.34"
$ws.Range("E27").Value = "'  -0.76%  "
$ws.Range("D28").Value = "'6.741"
$ws.Range("E28").Value = "'  -2.11%  "
$ws.Range("E29").Value = "'  -1.43%  "
$ws.Range("D30").Value = "'1.236"
$ws.Range("E30").Value = "'  -0.81%  "
$ws.Range("D31").Value = "'0.04881"
$ws.Range("E31").Value = "'  +0.29%  "
$ws.Range("D32").Value = "'3.302"
$ws.Range("E32").Value = "'  -0.40%  "
$ws.Range("D33").Value = "'3.187"
$ws.Range("E33").Value = "'  -1.81%  "
$ws.Range("D34").Value = "'1.559"
$ws.Range("E34").Value = "'  +0.71%  "
$ws.Range("D35").Value = "'2.374"
$ws.Range("E35").Value = "'  -0.33%  "
$ws.Range("D36").Value = "'0.8947"
$ws.Range("E36").Value = "'  -2.87%  "
$ws.Range("D37").Value = "'1.131.62"
$ws.Range("E37").Value = "'  +1.07%  "
$ws.Range("D38").Value = "'0.5494"
$ws.Range("E38").Value = "'  -1.87%  "
$ws.Range("D39").Value = "'2.531"
$ws.Range("E39").Value = "'  -1.56%  "
$ws.Range("D40").Value = "'0.01562"
$ws.Range("E40").Value = "'  -1.08%  "
$ws.Range("D41").Value = "'1.003"
$ws.Range("E41").Value = "'  +0.09%  "
$ws.Range("D42").Value = "'5.589"
$ws.Range("E42").Value = "'  +0.08%  "
$ws.Range("D43").Value = "'0.7946"
$ws.Range("E43").Value = "'  -1.94%  "
$ws.Range("D44").Value = "'97.33"
$ws.Range("D45").Value = "'1.772.68"
$ws.Range("E45").Value = "'  -0.57%  "
$ws.Range("D46").Value = "'0.0₈111"
$ws.Range("E46").Value = "'  -7.73%  "
$ws.Range("D47").Value = "'0.4422"
$ws.Range("E47").Value = "'  -2.51%  "
$ws.Range("D48").Value = "'54.81"
$ws.Range("E48").Value = "'  -1.17%  "
$ws.Range("D49").Value = "'0.05066"
$ws.Range("E49").Value = "'  -3.01%  "
$ws.Range("D50").Value = "'7.581"
$ws.Range("E50").Value = "'  +1.46%  "
$ws.Range("E51").Value = "'  -0.33%  "
